# Delete string resources and update file properties
#
# The English "differentiation algorithm" entries (strRadBackwardOne,
# strRadCentralFive, strRadCentralThree, strRadForwardOne + their English
# text, plus the now-orphaned comment "No need, since they are already
# listed in strDifferentiationAlgorithms") were removed from the
# translation table. In the worksheet this corresponds to deleting the
# four rows that held those entries; Excel prunes the now-unused shared
# strings and resizes the bound table/dimension automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 132 (strRadBackwardOne), 133 (strRadCentralFive), 134 (strRadCentralThree)
# are contiguous; row 136 (strRadForwardOne) follows row 135 (strRadCurrentCulture)
# which is kept. Delete the contiguous block first, then the single row -
# after the first delete, the old row 136 has already shifted up to row 133.
$ws.Rows("132:134").Delete() | Out-Null
$ws.Rows("133").Delete() | Out-Null

# Column B ("Key") was widened slightly in the source file.
$ws.Columns("B").ColumnWidth = 32.666666666666664
